$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.409.56"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "1.785.57"
$ws.Range("E3").Value = "  +4.07%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'313.71"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.5379"
$ws.Range("E7").Value = "  +12.83%  "
$ws.Range("D8").Value = "'0.3782"
$ws.Range("E8").Value = "  +9.02%  "
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("D10").Value = "'0.07429"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").Value = "'1.101"
$ws.Range("E11").Value = "  +5.58%  "
$ws.Range("D12").Value = "'0.9994"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'20.77"
$ws.Range("E13").Value = "  +4.87%  "
$ws.Range("D14").Value = "'6.128"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("D15").Value = "1.779.03"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("D16").Value = "'7.015"
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("D17").Value = "'89.81"
$ws.Range("E17").Value = "  +3.95%  "
$ws.Range("D18").Value = "'0.00001059"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").Value = "'0.06445"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "'16.86"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").Value = "'5.924"
$ws.Range("E22").Value = "  +5.49%  "
$ws.Range("D23").Value = "27.440.40"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").Value = "'11.22"
$ws.Range("E24").Value = "  +4.54%  "
$ws.Range("D25").Value = "'2.088"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'156.07"
$ws.Range("E26").Value = "  +3.23%  "
$ws.Range("D27").Value = "'20.26"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "'2.377"
$ws.Range("E28").Value = "  +13.61%  "
$ws.Range("D29").Value = "1.982.23"
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("D30").Value = "'121.29"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "'1.086"
$ws.Range("E31").Value = "  +5.64%  "
$ws.Range("D32").Value = "'0.1030"
$ws.Range("E32").Value = "  +12.63%  "
$ws.Range("D33").Value = "'5.614"
$ws.Range("E33").Value = "  +5.63%  "
$ws.Range("D34").Value = "'3.625"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'0.02262"
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("D36").Value = "'0.05992"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'11.33"
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2065"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.930"
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("D40").Value = "'8.341"
$ws.Range("E40").Value = "  +12.18%  "
$ws.Range("D41").Value = "'0.6143"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").Value = "'1.424"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'0.9988"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.138"
$ws.Range("E44").Value = "  +5.01%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.28"
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5796"
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.630"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'121.23"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.899"
$ws.Range("E49").Value = "  +3.89%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.130"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06728"
$ws.Range("E51").Value = "  +1.13%  "

Write-Output "Applied 123 cell updates"
